# Insert a new column "eln_id" between "sample_batch" (B) and "test_id" (old C,
# now D), matching the template's new constraint column, and move the active
# selection as it was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column B's width before the insert shifts everything, so the new
# column C can inherit it (Excel's "insert column" normally carries the
# left-hand column's formatting onto the freshly inserted one).
$bWidth = $ws.Columns("B").ColumnWidth

# Insert a new blank column at C; everything from the old C onward (test_id,
# study_type, target, ... data_source) shifts one column to the right.
$ws.Columns("C:C").Insert()

# New column C takes on column B's width.
$ws.Columns("C").ColumnWidth = $bWidth

# Header for the newly inserted column.
$ws.Range("C1").Value = "eln_id"

# Selection left where the author's session ended up.
$ws.Range("E9").Select()
